# "Fruta / hortaliza, semanal" -- the weekly refresh re-derives each
# listing's per-record facts (date, quality, volume, prices, unit,
# origin, $/kg, kg/unit) and re-distributes them across the existing
# data rows; the identifying columns (market/region/product/category/
# variety, A:C & E:K) stay put. Net effect on this sheet: rows 2-14
# keep their row position but columns D and L:T get reshuffled to the
# corrected record order.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The columns that move as a unit for each record.
$cols = @("D", "L", "M", "N", "O", "P", "Q", "R", "S", "T")

$firstRow = 2
$lastRow = 14

# 1) Snapshot every data row's movable columns BEFORE any writes, so
#    reads are never contaminated by an earlier write in this pass.
$snapshot = @{}
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $rowVals = @()
    foreach ($c in $cols) {
        $rowVals += $ws.Range("$c$r").Value2
    }
    $snapshot[$r] = $rowVals
}

# 2) Target row -> source row whose snapshot it should now receive.
$rowMap = @{
    2  = 3
    3  = 11
    4  = 8
    5  = 10
    6  = 14
    7  = 2
    8  = 6
    9  = 12
    10 = 9
    11 = 13
    12 = 4
    13 = 5
    14 = 7
}

# 3) Write each target row's cells from its mapped source snapshot.
foreach ($targetRow in ($rowMap.Keys | Sort-Object)) {
    $sourceRow = $rowMap[$targetRow]
    $values = $snapshot[$sourceRow]
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $ws.Range("$($cols[$i])$targetRow").Value = $values[$i]
    }
}

"Reshuffled rows $firstRow-$lastRow across columns $($cols -join ', ')"
